# Updated cryptos list on Tue Aug  8 09:20:57 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.155.12'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '1.829.34'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.92'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6166'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07341'
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2907'
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.18'
$ws.Range("E10").Value = '  +0.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07634'
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("D12").Value = '1.832.82'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.973'
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("E14").Value = '  -0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.28'
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008973'
$ws.Range("E16").Value = '  -2.05%  '
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").Value = '29.148.74'
$ws.Range("E18").Value = '  +0.39%  '
$ws.Range("D19").Value = '2.078.63'
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '235.55'
$ws.Range("E20").Value = '  +1.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.46'
$ws.Range("E21").Value = '  -1.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.349'
$ws.Range("E23").Value = '  +2.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.60'
$ws.Range("E25").Value = '  -0.30%  '
$ws.Range("E26").Value = '  -1.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.511'
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.59'
$ws.Range("E28").Value = '  -1.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.488'
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05857'
$ws.Range("E30").Value = '  +5.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.223'
$ws.Range("E31").Value = '  +1.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.080'
$ws.Range("E32").Value = '  -0.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.080'
$ws.Range("E33").Value = '  -1.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.858'
$ws.Range("E34").Value = '  +1.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7250'
$ws.Range("E35").Value = '  -1.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.136'
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.610'
$ws.Range("E37").Value = '  -1.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.857'
$ws.Range("E38").Value = '  +3.14%  '
$ws.Range("D39").Value = '1.229.49'
$ws.Range("E39").Value = '  +1.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01759'
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.199'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9032'
$ws.Range("E42").Value = '  +1.55%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.75'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = '1.983.00'
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.69'
$ws.Range("E46").Value = '  +0.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5040'
$ws.Range("E47").Value = '  -0.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4050'
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.182'
$ws.Range("E49").Value = '  +0.76%  '
$ws.Range("E50").Value = '  -3.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1135'
$ws.Range("E51").Value = '  +2.99%  '
